$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "Förändrad" (C) column for every existing data row (2..411):
#    45192 -> 45202
$ws.Range("C2:C411").Value = 45202

# 2) Row 411 picks up an explicit row height (ht="15" customHeight="1")
$ws.Rows.Item(411).RowHeight = 15

# 3) Append five new report rows (412..416) with the same column layout as
#    the existing rows: A..E, G..Q populated, R left blank (wrap-text style).
$newRows = @(
    @{ Row=412; A="A 46026-2023"; B=45196; C=45202; G=0.9 },
    @{ Row=413; A="A 46236-2023"; B=45196; C=45202; G=8.300000000000001 },
    @{ Row=414; A="A 46395-2023"; B=45197; C=45202; G=2 },
    @{ Row=415; A="A 46421-2023"; B=45197; C=45202; G=2.3 },
    @{ Row=416; A="A 46877-2023"; B=45200; C=45202; G=8.800000000000001 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.A

    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 2).NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 3).NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($row, 4).Value = "HALLANDS LÄN"
    $ws.Cells.Item($row, 5).Value = "VARBERG"

    $ws.Cells.Item($row, 7).Value = $r.G

    for ($col = 8; $col -le 17; $col++) {
        $ws.Cells.Item($row, $col).Value = 0
    }

    $ws.Cells.Item($row, 18).WrapText = $true

    if ($row -le 415) {
        $ws.Rows.Item($row).RowHeight = 15
    }
}
